$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-14 05:01:06"
$wsZh.Range("G2").Value = "2016-01-14 05:02:37"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-14 05:01:31"
$wsDe.Range("G2").Value = "2016-01-14 05:03:15"
